# Fruta / hortaliza, semanal
# Insert 3 new weekly report rows (Calibre 100, 2023-07-20) above the existing
# "Región de O'Higgins" block that starts at row 223, pushing the rest of the
# table down by three rows (223-258 -> 226-261).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at row 223 (each insert pushes the following rows down).
$ws.Rows.Item(223).Insert()
$ws.Rows.Item(224).Insert()
$ws.Rows.Item(225).Insert()

# Common (repeated) column values for this market / product.
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$prodId    = 100104
$producto  = "Frutos de pepita"
$catId     = 100104002
$categoria = "Manzana"
$unidad    = "$/caja 18 kilos embalada"
$origen    = "Región de O'Higgins"
$kgUnidad  = 18
$fecha     = 45127

$filas = @(
    @{ Row = 223; Variedad = "Fuji royal";   Calidad = "Calibre 100"; Volumen = 270; PMin = 22000; PMax = 23000; PProm = 22500; PKg = 1250 },
    @{ Row = 224; Variedad = "Granny Smith"; Calidad = "Calibre 100"; Volumen = 270; PMin = 22000; PMax = 23000; PProm = 22500; PKg = 1250 },
    @{ Row = 225; Variedad = "Royal Gala";   Calidad = "Calibre 100"; Volumen = 300; PMin = 22000; PMax = 23000; PProm = 22500; PKg = 1250 }
)

foreach ($fila in $filas) {
    $r = $fila.Row
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $prodId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $catId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $fila.Variedad
    $ws.Cells.Item($r, 12).Value = $fila.Calidad
    $ws.Cells.Item($r, 13).Value = $fila.Volumen
    $ws.Cells.Item($r, 14).Value = $fila.PMin
    $ws.Cells.Item($r, 15).Value = $fila.PMax
    $ws.Cells.Item($r, 16).Value = $fila.PProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $fila.PKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
